# Applies the "cleaning up for opensource" edit:
#  - Row heights for data rows 2-19 change from 18 to 18.75
#  - Row heights for data rows 20-22 change from 18 to 19.5
#  - The data-table font (used by every cell except the header label cells
#    B1:R1) is recolored from theme color 1 to explicit RGB black
#    (cosmetically identical, but stored as an explicit rgb value rather
#    than a theme reference)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights -----------------------------------------------------
for ($r = 2; $r -le 19; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
for ($r = 20; $r -le 22; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}

# --- Font color cleanup ------------------------------------------------
# Every cell in the table (A1:R52, header row included) ends up using the
# font that carried <color theme="1"/>; it is switched to an explicit
# black RGB color instead of the theme reference.
$ws.Range("A1:R52").Font.Color = 0
